$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3452351689338684
$ws.Range("B1").Value = 0.3323192894458771
$ws.Range("C1").Value = 0.3397984802722931
$ws.Range("D1").Value = 0.4460178315639496
$ws.Range("E1").Value = 0.6353998780250549
